$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 9453
$ws.Range("D2").Value = 8379
$ws.Range("E2").Value = 0.8863852745160267
$ws.Range("F2").Value = 0.8840472673559823
$ws.Range("G2").Value = 0.09675112722510387
$ws.Range("H2").Value = 0.08553256963696407
$ws.Range("I2").Value = 41181848.30283703
$ws.Range("J2").Value = 14396092.09807451
$ws.Range("L2").Value = 14396092.09807451
$ws.Range("M2").Value = 55577940.40091153
$ws.Range("N2").Value = 800091928.4872
$ws.Range("O2").Value = 782392121.4832001
$ws.Range("P2").Value = 0.01799304753054364
$ws.Range("Q2").Value = 0.01840009849636968

# Row 3
$ws.Range("C3").Value = 9644
$ws.Range("D3").Value = 8563
$ws.Range("E3").Value = 0.8879095810866861
$ws.Range("F3").Value = 0.8858886819780675
$ws.Range("G3").Value = 0.09538540115674442
$ws.Range("H3").Value = 0.08450084731069751
$ws.Range("I3").Value = 43122511.22151443
$ws.Range("J3").Value = 15096506.55049641
$ws.Range("L3").Value = 15096506.55049641
$ws.Range("M3").Value = 58219017.77201083
$ws.Range("N3").Value = 837860675.346328
$ws.Range("O3").Value = 820380499.322258
$ws.Range("P3").Value = 0.01801791991760002
$ws.Range("Q3").Value = 0.01840183495703287

# Row 4
$ws.Range("C4").Value = 9824
$ws.Range("D4").Value = 8726
$ws.Range("E4").Value = 0.8882328990228013
$ws.Range("F4").Value = 0.8851694055589369
$ws.Range("G4").Value = 0.09438221778644625
$ws.Range("H4").Value = 0.08354425161336274
$ws.Range("I4").Value = 45115835.6214844
$ws.Range("J4").Value = 15773246.61981758
$ws.Range("L4").Value = 15773246.61981758
$ws.Range("M4").Value = 60889082.24130198
$ws.Range("N4").Value = 874134762.184269
$ws.Range("O4").Value = 856685814.1783152
$ws.Range("P4").Value = 0.01804441065860799
$ws.Range("Q4").Value = 0.01841193861129403

# Row 5
$ws.Range("D5").Value = 8915
$ws.Range("E5").Value = 0.8889221258350782
$ws.Range("F5").Value = 0.8867117565148199
$ws.Range("G5").Value = 0.09325617978731912
$ws.Range("H5").Value = 0.08269135098507559
$ws.Range("I5").Value = 47276888.58374348
$ws.Range("J5").Value = 16520699.39275815
$ws.Range("L5").Value = 16520699.39275815
$ws.Range("M5").Value = 63797587.97650164
$ws.Range("N5").Value = 914424195.1217525
$ws.Range("O5").Value = 896938089.6587793
$ws.Range("P5").Value = 0.01806677850486937
$ws.Range("Q5").Value = 0.01841899634237085

# Row 6
$ws.Range("C6").Value = 10229
$ws.Range("D6").Value = 9106
$ws.Range("E6").Value = 0.8902140971746993
$ws.Range("F6").Value = 0.8880436902672127
$ws.Range("G6").Value = 0.09212643515120535
$ws.Range("H6").Value = 0.08181229944283946
$ws.Range("I6").Value = 49566607.15969561
$ws.Range("J6").Value = 17294514.02090722
$ws.Range("L6").Value = 17294514.02090722
$ws.Range("M6").Value = 66861121.18060283
$ws.Range("N6").Value = 955116216.3787864
$ws.Range("O6").Value = 937524390.505605
$ws.Range("P6").Value = 0.01810723524983942
$ws.Range("Q6").Value = 0.01844700169515624
